$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4193.6665
$ws.Range("I88").Value = 4193.6665
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 4193.6665
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3787.6665
$ws.Range("H91").Value = 4193.6665
$ws.Range("I91").Value = 4193.6665
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 4193.6665
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -2789.6665
$ws.Range("H103").Value = 1625.6471
$ws.Range("I103").Value = 346.3
$ws.Range("J103").Value = 3453.2856
$ws.Range("K103").Value = 1038.9
$ws.Range("L103").Value = 10359.8568
$ws.Range("M103").Value = -452.9000000000001
$ws.Range("N103").Value = -11531.8568
$ws.Range("H116").Value = 7613.343
$ws.Range("I116").Value = 7325.6924
$ws.Range("J116").Value = 8444.333000000001
$ws.Range("K116").Value = 7325.6924
$ws.Range("L116").Value = 8444.333000000001
$ws.Range("M116").Value = -3883.6924
$ws.Range("N116").Value = -15328.333

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 189.5
$ws.Range("I38").Value = 189.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 189.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 277.5
$ws.Range("H61").Value = 3701.0847
$ws.Range("I61").Value = 2792.5862
$ws.Range("J61").Value = 4579.3
$ws.Range("K61").Value = 2792.5862
$ws.Range("L61").Value = 4579.3
$ws.Range("M61").Value = -2580.5862
$ws.Range("N61").Value = -5003.3
$ws.Range("H74").Value = 1766865.9
$ws.Range("I74").Value = 2273747.2
$ws.Range("J74").Value = 77261.25
$ws.Range("K74").Value = 2273747.2
$ws.Range("L74").Value = 77261.25
$ws.Range("M74").Value = -2272873.2
$ws.Range("N74").Value = -79009.25
$ws.Range("H77").Value = 1766865.9
$ws.Range("I77").Value = 2273747.2
$ws.Range("J77").Value = 77261.25
$ws.Range("K77").Value = 11368736
$ws.Range("L77").Value = 386306.25
$ws.Range("M77").Value = -11364368
$ws.Range("N77").Value = -395042.25
$ws.Range("H102").Value = 58827412
$ws.Range("I102").Value = 66670510
$ws.Range("J102").Value = 4149
$ws.Range("K102").Value = 66670510
$ws.Range("L102").Value = 4149
$ws.Range("M102").Value = -66668888
$ws.Range("N102").Value = -7393
$ws.Range("H132").Value = 20466.2
$ws.Range("I132").Value = 21713.785
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 65141.355
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -62611.355
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 3701.0847
$ws.Range("I136").Value = 2792.5862
$ws.Range("J136").Value = 4579.3
$ws.Range("K136").Value = 8377.758600000001
$ws.Range("L136").Value = 13737.9
$ws.Range("M136").Value = -5827.758600000001
$ws.Range("N136").Value = -18837.9

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2939.5
$ws.Range("I20").Value = 2939.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2939.5
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -2692.5
$ws.Range("H86").Value = 22729658
$ws.Range("I86").Value = 34485196
$ws.Range("J86").Value = 2285.6
$ws.Range("K86").Value = 34485196
$ws.Range("L86").Value = 2285.6
$ws.Range("M86").Value = -34484073
$ws.Range("N86").Value = -4531.6
$ws.Range("H89").Value = 22729658
$ws.Range("I89").Value = 34485196
$ws.Range("J89").Value = 2285.6
$ws.Range("K89").Value = 172425980
$ws.Range("L89").Value = 11428
$ws.Range("M89").Value = -172420364
$ws.Range("N89").Value = -22660
$ws.Range("H107").Value = 12299.038
$ws.Range("I107").Value = 15501.75
$ws.Range("J107").Value = 1623.3334
$ws.Range("K107").Value = 15501.75
$ws.Range("L107").Value = 1623.3334
$ws.Range("M107").Value = -13581.75
$ws.Range("N107").Value = -5463.3334
$ws.Range("H134").Value = 2818.08
$ws.Range("I134").Value = 2783.5789
$ws.Range("J134").Value = 2927.3333
$ws.Range("K134").Value = 8350.736699999999
$ws.Range("L134").Value = 8781.999899999999
$ws.Range("M134").Value = -5815.736699999999
$ws.Range("N134").Value = -13851.9999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1919.091
$ws.Range("I22").Value = 1477.75
$ws.Range("J22").Value = 2171.2856
$ws.Range("K22").Value = 1477.75
$ws.Range("L22").Value = 2171.2856
$ws.Range("M22").Value = -1127.75
$ws.Range("N22").Value = -2871.2856
$ws.Range("H58").Value = 2902.4
$ws.Range("I58").Value = 2535.0625
$ws.Range("J58").Value = 3555.4443
$ws.Range("K58").Value = 2535.0625
$ws.Range("L58").Value = 3555.4443
$ws.Range("M58").Value = -2332.0625
$ws.Range("N58").Value = -3961.4443
$ws.Range("H99").Value = 5045.32
$ws.Range("I99").Value = 4843.6665
$ws.Range("J99").Value = 5347.8
$ws.Range("K99").Value = 4843.6665
$ws.Range("L99").Value = 5347.8
$ws.Range("M99").Value = -3345.6665
$ws.Range("N99").Value = -8343.799999999999
$ws.Range("H107").Value = 29439926
$ws.Range("I107").Value = 35747670
$ws.Range("J107").Value = 3778.8333
$ws.Range("K107").Value = 35747670
$ws.Range("L107").Value = 3778.8333
$ws.Range("M107").Value = -35745750
$ws.Range("N107").Value = -7618.8333
$ws.Range("H122").Value = 4124.75
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 4666.3335
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 13999.0005
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -18899.0005
$ws.Range("H126").Value = 5045.32
$ws.Range("I126").Value = 4843.6665
$ws.Range("J126").Value = 5347.8
$ws.Range("K126").Value = 14530.9995
$ws.Range("L126").Value = 16043.4
$ws.Range("M126").Value = -12060.9995
$ws.Range("N126").Value = -20983.4
$ws.Range("H132").Value = 4462.9473
$ws.Range("I132").Value = 3424.75
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 10274.25
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -7744.25
$ws.Range("N132").Value = -35060
$ws.Range("H134").Value = 4134.1353
$ws.Range("I134").Value = 4334.7
$ws.Range("J134").Value = 3274.5715
$ws.Range("K134").Value = 13004.1
$ws.Range("L134").Value = 9823.7145
$ws.Range("M134").Value = -10469.1
$ws.Range("N134").Value = -14893.7145
$ws.Range("H136").Value = 2902.4
$ws.Range("I136").Value = 2535.0625
$ws.Range("J136").Value = 3555.4443
$ws.Range("K136").Value = 7605.1875
$ws.Range("L136").Value = 10666.3329
$ws.Range("M136").Value = -5055.1875
$ws.Range("N136").Value = -15766.3329

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2825.7144
$ws.Range("I131").Value = 1328.7693
$ws.Range("J131").Value = 5258.25
$ws.Range("K131").Value = 3986.3079
$ws.Range("L131").Value = 15774.75
$ws.Range("M131").Value = 1053.6921
$ws.Range("N131").Value = -25854.75
$ws.Range("H134").Value = 3830.5
$ws.Range("I134").Value = 3830.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11491.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6421.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 20933
$ws.Range("I31").Value = 20933
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 20933
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -20641
$ws.Range("H37").Value = 20933
$ws.Range("I37").Value = 20933
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 20933
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -20656
$ws.Range("H46").Value = 22689.5
$ws.Range("I46").Value = 6047
$ws.Range("J46").Value = 39332
$ws.Range("K46").Value = 6047
$ws.Range("L46").Value = 39332
$ws.Range("M46").Value = -5891
$ws.Range("N46").Value = -39644
$ws.Range("H80").Value = 5847.55
$ws.Range("I80").Value = 2340.5833
$ws.Range("J80").Value = 11108
$ws.Range("K80").Value = 2340.5833
$ws.Range("L80").Value = 11108
$ws.Range("M80").Value = -1342.5833
$ws.Range("N80").Value = -13104
$ws.Range("H83").Value = 5847.55
$ws.Range("I83").Value = 2340.5833
$ws.Range("J83").Value = 11108
$ws.Range("K83").Value = 11702.9165
$ws.Range("L83").Value = 55540
$ws.Range("M83").Value = -6710.916499999999
$ws.Range("N83").Value = -65524
$ws.Range("H102").Value = 65479.61
$ws.Range("I102").Value = 129555.875
$ws.Range("J102").Value = 14218.6
$ws.Range("K102").Value = 129555.875
$ws.Range("L102").Value = 14218.6
$ws.Range("M102").Value = -127933.875
$ws.Range("N102").Value = -17462.6
$ws.Range("H126").Value = 48018.1
$ws.Range("I126").Value = 55955.59
$ws.Range("J126").Value = 3039
$ws.Range("K126").Value = 167866.77
$ws.Range("L126").Value = 9117
$ws.Range("M126").Value = -165396.77
$ws.Range("N126").Value = -14057
$ws.Range("H132").Value = 4693.231
$ws.Range("I132").Value = 4908.909
$ws.Range("J132").Value = 3507
$ws.Range("K132").Value = 14726.727
$ws.Range("L132").Value = 10521
$ws.Range("M132").Value = -12196.727
$ws.Range("N132").Value = -15581

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 11713.286
$ws.Range("I32").Value = 8999
$ws.Range("J32").Value = 13749
$ws.Range("K32").Value = 8999
$ws.Range("L32").Value = 13749
$ws.Range("M32").Value = -8682
$ws.Range("N32").Value = -14383
$ws.Range("H46").Value = 3473.4167
$ws.Range("I46").Value = 867.25
$ws.Range("J46").Value = 3799.1875
$ws.Range("K46").Value = 867.25
$ws.Range("L46").Value = 3799.1875
$ws.Range("M46").Value = -679.25
$ws.Range("N46").Value = -4175.1875
$ws.Range("H61").Value = 17858416
$ws.Range("I61").Value = 22728314
$ws.Range("J61").Value = 2122.8333
$ws.Range("K61").Value = 22728314
$ws.Range("L61").Value = 2122.8333
$ws.Range("M61").Value = -22728112
$ws.Range("N61").Value = -2526.8333
$ws.Range("H113").Value = 17858416
$ws.Range("I113").Value = 22728314
$ws.Range("J113").Value = 2122.8333
$ws.Range("K113").Value = 22728314
$ws.Range("L113").Value = 2122.8333
$ws.Range("M113").Value = -22726144
$ws.Range("N113").Value = -6462.8333
$ws.Range("H122").Value = 3371.7036
$ws.Range("I122").Value = 3321.6
$ws.Range("J122").Value = 3998
$ws.Range("K122").Value = 9964.799999999999
$ws.Range("L122").Value = 11994
$ws.Range("M122").Value = -7514.799999999999
$ws.Range("N122").Value = -16894
$ws.Range("H132").Value = 2420.8333
$ws.Range("I132").Value = 2051.9167
$ws.Range("J132").Value = 4634.3335
$ws.Range("K132").Value = 6155.750100000001
$ws.Range("L132").Value = 13903.0005
$ws.Range("M132").Value = -3625.750100000001
$ws.Range("N132").Value = -18963.0005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 13669.25
$ws.Range("I74").Value = 14000
$ws.Range("J74").Value = 13559
$ws.Range("K74").Value = 14000
$ws.Range("L74").Value = 13559
$ws.Range("M74").Value = -13064
$ws.Range("N74").Value = -15431
$ws.Range("H77").Value = 13669.25
$ws.Range("I77").Value = 14000
$ws.Range("J77").Value = 13559
$ws.Range("K77").Value = 42000
$ws.Range("L77").Value = 40677
$ws.Range("M77").Value = -37320
$ws.Range("N77").Value = -50037
$ws.Range("H100").Value = 445.57144
$ws.Range("I100").Value = 358.92307
$ws.Range("J100").Value = 586.375
$ws.Range("K100").Value = 717.84614
$ws.Range("L100").Value = 1172.75
$ws.Range("M100").Value = -176.84614
$ws.Range("N100").Value = -2254.75

Write-Host "applied edits"
